# Update the "想去人数" (number of people interested) values in column F
# for rows 2-4 on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 25
    $ws.Range("F3").Value = 50
    $ws.Range("F4").Value = 11
}
